# Update "想去人数" (F column) figures across the sheets to reflect the
# latest generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(6, 6).Value  = 1237
$ws1.Cells.Item(9, 6).Value  = 964
$ws1.Cells.Item(11, 6).Value = 2328
$ws1.Cells.Item(12, 6).Value = 30
$ws1.Cells.Item(13, 6).Value = 1147
$ws1.Cells.Item(14, 6).Value = 879
$ws1.Cells.Item(16, 6).Value = 867
$ws1.Cells.Item(17, 6).Value = 1039
$ws1.Cells.Item(21, 6).Value = 726
$ws1.Cells.Item(22, 6).Value = 157
$ws1.Cells.Item(23, 6).Value = 426
$ws1.Cells.Item(24, 6).Value = 1074
$ws1.Cells.Item(26, 6).Value = 490
$ws1.Cells.Item(29, 6).Value = 275
$ws1.Cells.Item(31, 6).Value = 645
$ws1.Cells.Item(32, 6).Value = 2737
$ws1.Cells.Item(33, 6).Value = 443
$ws1.Cells.Item(38, 6).Value = 1531
$ws1.Cells.Item(40, 6).Value = 137
$ws1.Cells.Item(42, 6).Value = 120
$ws1.Cells.Item(47, 6).Value = 69

# ---- Sheet: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(11, 6).Value = 4409
$ws2.Cells.Item(13, 6).Value = 16
$ws2.Cells.Item(14, 6).Value = 128

# ---- Sheet: 本地生活 (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 2234
$ws3.Cells.Item(3, 6).Value = 688
$ws3.Cells.Item(4, 6).Value = 653

# ---- Sheet: 全部类型 (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value  = 2234
$ws4.Cells.Item(5, 6).Value  = 653
$ws4.Cells.Item(6, 6).Value  = 1237
$ws4.Cells.Item(10, 6).Value = 964
$ws4.Cells.Item(11, 6).Value = 2328
$ws4.Cells.Item(12, 6).Value = 30
$ws4.Cells.Item(13, 6).Value = 1147
$ws4.Cells.Item(14, 6).Value = 879
$ws4.Cells.Item(16, 6).Value = 867
$ws4.Cells.Item(17, 6).Value = 1039
$ws4.Cells.Item(24, 6).Value = 726
$ws4.Cells.Item(25, 6).Value = 157
$ws4.Cells.Item(26, 6).Value = 426
$ws4.Cells.Item(27, 6).Value = 1074
$ws4.Cells.Item(30, 6).Value = 490
$ws4.Cells.Item(33, 6).Value = 275
$ws4.Cells.Item(34, 6).Value = 2737
$ws4.Cells.Item(36, 6).Value = 443
$ws4.Cells.Item(38, 6).Value = 1531
$ws4.Cells.Item(40, 6).Value = 137
$ws4.Cells.Item(42, 6).Value = 16
$ws4.Cells.Item(43, 6).Value = 120
$ws4.Cells.Item(47, 6).Value = 69
